$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows 8 and 9 with CAMID / Error info
$ws.Range("B8").Value = "907E"
$ws.Range("D8").Value = "one day off in last record"

$ws.Range("B9").Value = "909A"
$ws.Range("D9").Value = "some should be 2020 not 2019 "

# Update the active selection to D9 as in the diff
$ws.Range("D9").Select()
